# Auto update Excel log
# Appends newly-logged sensor readings to the PIR, Humidity, Proximity and
# Camera sheets of the SeniorConnect master log workbook.

$wb = $excel.ActiveWorkbook

function Add-LogRows {
    param(
        [string]$SheetName,
        [int]$StartRow,
        [object[]]$Rows
    )

    $ws = $wb.Worksheets.Item($SheetName)

    for ($i = 0; $i -lt $Rows.Count; $i++) {
        $r = $StartRow + $i
        $row = $Rows[$i]

        # Prefix every value with an apostrophe so the COM layer stores each
        # one as literal text instead of auto-converting look-alike values
        # (dates, percentages, times, numbers) into typed cell values.
        for ($col = 1; $col -le 6; $col++) {
            $ws.Cells.Item($r, $col).Value = "'" + $row[$col - 1]
        }
    }
}

# ---------------------------------------------------------------------------
# PIR sheet: rows 141-153 (Bathroom / No Motion / Inactive)
# ---------------------------------------------------------------------------
$pirRows = @(
    @("2026-01-30", "16:18:33", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "16:18:35", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "16:18:40", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "16:18:45", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "16:18:50", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "16:18:55", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "16:19:00", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "16:19:05", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "16:19:11", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "16:19:15", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "16:19:20", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "16:19:25", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "16:19:30", "16:00", "Bathroom", "No Motion", "Inactive")
)
Add-LogRows "PIR" 141 $pirRows

# ---------------------------------------------------------------------------
# Humidity sheet: rows 108-114 (Bathroom / 87.3% / Active)
# ---------------------------------------------------------------------------
$humidityRows = @(
    @("2026-01-30", "16:18:33", "16:00", "Bathroom", "87.3%", "Active"),
    @("2026-01-30", "16:18:41", "16:00", "Bathroom", "87.3%", "Active"),
    @("2026-01-30", "16:18:50", "16:00", "Bathroom", "87.3%", "Active"),
    @("2026-01-30", "16:19:00", "16:00", "Bathroom", "87.3%", "Active"),
    @("2026-01-30", "16:19:11", "16:00", "Bathroom", "87.3%", "Active"),
    @("2026-01-30", "16:19:16", "16:00", "Bathroom", "87.3%", "Active"),
    @("2026-01-30", "16:19:31", "16:00", "Bathroom", "87.3%", "Active")
)
Add-LogRows "Humidity" 108 $humidityRows

# ---------------------------------------------------------------------------
# Proximity sheet: rows 42-47 (Living Room Main Door ENTER/EXIT events)
# ---------------------------------------------------------------------------
$proximityRows = @(
    @("2026-01-30", "16:18:34", "16:00", "Living Room Main Door", "EXIT", "User EXITED Living Room Main Door"),
    @("2026-01-30", "16:18:41", "16:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door"),
    @("2026-01-30", "16:18:56", "16:00", "Living Room Main Door", "EXIT", "User EXITED Living Room Main Door"),
    @("2026-01-30", "16:19:00", "16:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door"),
    @("2026-01-30", "16:19:10", "16:00", "Living Room Main Door", "EXIT", "User EXITED Living Room Main Door"),
    @("2026-01-30", "16:19:12", "16:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door")
)
Add-LogRows "Proximity" 42 $proximityRows

# ---------------------------------------------------------------------------
# Camera sheet: rows 42-47 (Living Room Main Door image captures)
# ---------------------------------------------------------------------------
$cameraRows = @(
    @("2026-01-30", "16:18:34", "16:00", "Living Room Main Door", "Image Captured (EXIT)", "Active"),
    @("2026-01-30", "16:18:41", "16:00", "Living Room Main Door", "Image Captured (ENTER)", "Active"),
    @("2026-01-30", "16:18:56", "16:00", "Living Room Main Door", "Image Captured (EXIT)", "Active"),
    @("2026-01-30", "16:18:59", "16:00", "Living Room Main Door", "Image Captured (ENTER)", "Active"),
    @("2026-01-30", "16:19:10", "16:00", "Living Room Main Door", "Image Captured (EXIT)", "Active"),
    @("2026-01-30", "16:19:12", "16:00", "Living Room Main Door", "Image Captured (ENTER)", "Active")
)
Add-LogRows "Camera" 42 $cameraRows
